# Scheduled runner refresh: re-pull current market-board prices and
# re-derive the dependent Leve profit columns (H-N) per row, across all
# eight crafting-job sheets. Only numeric market/profit cells change;
# item/leve metadata columns (A-G) are untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H19").Value = 2197.3635
$ws.Range("I19").Value = 2600
$ws.Range("K19").Value = 2600
$ws.Range("M19").Value = -2425

$ws.Range("H41").Value = 875
$ws.Range("I41").Value = 875
$ws.Range("K41").Value = 875
$ws.Range("M41").Value = -435

$ws.Range("H70").Value = 1654.4445
$ws.Range("I70").Value = 1498.3334
$ws.Range("J70").Value = 1732.5
$ws.Range("K70").Value = 4495.0002
$ws.Range("L70").Value = 5197.5
$ws.Range("M70").Value = -4225.0002
$ws.Range("N70").Value = -5737.5

$ws.Range("H73").Value = 1654.4445
$ws.Range("I73").Value = 1498.3334
$ws.Range("J73").Value = 1732.5
$ws.Range("K73").Value = 4495.0002
$ws.Range("L73").Value = 5197.5
$ws.Range("M73").Value = -3559.0002
$ws.Range("N73").Value = -7069.5

$ws.Range("H80").Value = 1658.4445
$ws.Range("I80").Value = 499.2
$ws.Range("J80").Value = 2104.3076
$ws.Range("K80").Value = 1497.6
$ws.Range("L80").Value = 6312.9228
$ws.Range("M80").Value = -499.5999999999999
$ws.Range("N80").Value = -8308.9228

$ws.Range("H82").Value = 176.66667
$ws.Range("I82").Value = 176.66667
$ws.Range("K82").Value = 530.00001
$ws.Range("M82").Value = -124.00001

$ws.Range("H83").Value = 1658.4445
$ws.Range("I83").Value = 499.2
$ws.Range("J83").Value = 2104.3076
$ws.Range("K83").Value = 4492.8
$ws.Range("L83").Value = 18938.7684
$ws.Range("M83").Value = 499.1999999999998
$ws.Range("N83").Value = -28922.7684

$ws.Range("H85").Value = 176.66667
$ws.Range("I85").Value = 176.66667
$ws.Range("K85").Value = 530.00001
$ws.Range("M85").Value = 873.99999

$ws.Range("H100").Value = 2334.1667
$ws.Range("I100").Value = 1666.6666
$ws.Range("J100").Value = 3001.6667
$ws.Range("K100").Value = 1666.6666
$ws.Range("L100").Value = 3001.6667
$ws.Range("M100").Value = -1125.6666
$ws.Range("N100").Value = -4083.6667

$ws.Range("H116").Value = 6915.75
$ws.Range("J116").Value = 6999
$ws.Range("L116").Value = 6999
$ws.Range("N116").Value = -13883

$ws.Range("H132").Value = 1957.5834
$ws.Range("I132").Value = 1957.5834
$ws.Range("K132").Value = 5872.7502
$ws.Range("M132").Value = -3342.7502

$ws.Range("H137").Value = 1499.6666
$ws.Range("I137").Value = 1249.5
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 3748.5
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -1198.5
$ws.Range("N137").Value = -11100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H61").Value = 1204
$ws.Range("I61").Value = 1204
$ws.Range("K61").Value = 1204
$ws.Range("M61").Value = -992

$ws.Range("H97").Value = 841.8889
$ws.Range("I97").Value = 788.1667
$ws.Range("J97").Value = 949.3333
$ws.Range("K97").Value = 788.1667
$ws.Range("L97").Value = 949.3333
$ws.Range("M97").Value = -292.1667
$ws.Range("N97").Value = -1941.3333

$ws.Range("H102").Value = 1186.375
$ws.Range("I102").Value = 1186.375
$ws.Range("K102").Value = 1186.375
$ws.Range("M102").Value = 435.625

$ws.Range("H132").Value = 5065.5
$ws.Range("I132").Value = 3754
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 11262
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -8732
$ws.Range("N132").Value = -32060

$ws.Range("H136").Value = 1204
$ws.Range("I136").Value = 1204
$ws.Range("K136").Value = 3612
$ws.Range("M136").Value = -1062

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 766.1667
$ws.Range("I94").Value = 780.875
$ws.Range("J94").Value = 736.75
$ws.Range("K94").Value = 780.875
$ws.Range("L94").Value = 736.75
$ws.Range("M94").Value = -329.875
$ws.Range("N94").Value = -1638.75

$ws.Range("H99").Value = 1485.4615
$ws.Range("I99").Value = 1456.4546
$ws.Range("J99").Value = 1645
$ws.Range("K99").Value = 1456.4546
$ws.Range("L99").Value = 1645
$ws.Range("M99").Value = 41.54539999999997
$ws.Range("N99").Value = -4641

$ws.Range("H105").Value = 2477
$ws.Range("I105").Value = 2461.625
$ws.Range("K105").Value = 2461.625
$ws.Range("M105").Value = -714.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4426.6
$ws.Range("I16").Value = 3110
$ws.Range("K16").Value = 3110
$ws.Range("M16").Value = -2823

$ws.Range("H113").Value = 4426.6
$ws.Range("I113").Value = 3110
$ws.Range("K113").Value = 3110
$ws.Range("M113").Value = -940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 2495
$ws.Range("J32").Value = 2495
$ws.Range("L32").Value = 7485
$ws.Range("N32").Value = -8051

$ws.Range("H33").Value = 713.3333
$ws.Range("I33").Value = 41
$ws.Range("K33").Value = 246
$ws.Range("M33").Value = 37

$ws.Range("H40").Value = 67.2
$ws.Range("I40").Value = 72.5
$ws.Range("J40").Value = 46
$ws.Range("K40").Value = 290
$ws.Range("L40").Value = 184
$ws.Range("M40").Value = -221
$ws.Range("N40").Value = -322

$ws.Range("H69").Value = 2000
$ws.Range("J69").Value = 2000
$ws.Range("L69").Value = 6000
$ws.Range("N69").Value = -7622

$ws.Range("H72").Value = 2000
$ws.Range("J72").Value = 2000
$ws.Range("L72").Value = 18000
$ws.Range("N72").Value = -26112

$ws.Range("H92").Value = 426.4
$ws.Range("I92").Value = 700
$ws.Range("J92").Value = 358
$ws.Range("K92").Value = 2100
$ws.Range("L92").Value = 1074
$ws.Range("M92").Value = -852
$ws.Range("N92").Value = -3570

$ws.Range("H98").Value = 3996.5
$ws.Range("J98").Value = 3996.5
$ws.Range("L98").Value = 11989.5
$ws.Range("N98").Value = -14985.5

$ws.Range("H103").Value = 918.75
$ws.Range("J103").Value = 1445.6
$ws.Range("L103").Value = 4336.799999999999
$ws.Range("N103").Value = -6094.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5217.5713
$ws.Range("I70").Value = 5126.5
$ws.Range("K70").Value = 5126.5
$ws.Range("M70").Value = -4856.5

$ws.Range("H73").Value = 5217.5713
$ws.Range("I73").Value = 5126.5
$ws.Range("K73").Value = 5126.5
$ws.Range("M73").Value = -4190.5

$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()

$ws.Range("H107").Value = 3317.375
$ws.Range("I107").Value = 1495
$ws.Range("J107").Value = 3577.7144
$ws.Range("K107").Value = 1495
$ws.Range("L107").Value = 3577.7144
$ws.Range("M107").Value = 425
$ws.Range("N107").Value = -7417.7144

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H136").Value = 50162.5
$ws.Range("J136").Value = 50162.5
$ws.Range("L136").Value = 150487.5
$ws.Range("N136").Value = -155587.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 781.7143
$ws.Range("I22").Value = 894
$ws.Range("J22").Value = 501
$ws.Range("K22").Value = 894
$ws.Range("L22").Value = 501
$ws.Range("M22").Value = -599
$ws.Range("N22").Value = -1091

$ws.Range("H27").Value = 781.7143
$ws.Range("I27").Value = 894
$ws.Range("J27").Value = 501
$ws.Range("K27").Value = 894
$ws.Range("L27").Value = 501
$ws.Range("M27").Value = -787
$ws.Range("N27").Value = -715

$ws.Range("H68").Value = 3648.5

$ws.Range("H71").Value = 3648.5

$ws.Range("H93").Value = 664.8333
$ws.Range("I93").Value = 698
$ws.Range("K93").Value = 698
$ws.Range("M93").Value = 550

$ws.Range("H100").Value = 3700.2
$ws.Range("I100").Value = 3100.4285
$ws.Range("J100").Value = 5099.6665
$ws.Range("K100").Value = 3100.4285
$ws.Range("L100").Value = 5099.6665
$ws.Range("M100").Value = -2559.4285
$ws.Range("N100").Value = -6181.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 44750.25
$ws.Range("J62").Value = 36333
$ws.Range("L62").Value = 36333
$ws.Range("N62").Value = -37581

$ws.Range("H65").Value = 44750.25
$ws.Range("J65").Value = 36333
$ws.Range("L65").Value = 181665
$ws.Range("N65").Value = -187905

$ws.Range("H81").Value = 1252684.1
$ws.Range("I81").Value = 1245
$ws.Range("K81").Value = 2490
$ws.Range("M81").Value = -1429

$ws.Range("H84").Value = 1252684.1
$ws.Range("I84").Value = 1245
$ws.Range("K84").Value = 12450
$ws.Range("M84").Value = -7146

$ws.Range("H96").Value = 1700
$ws.Range("J96").Value = 1700
$ws.Range("L96").Value = 1700
$ws.Range("N96").Value = -4446

$ws.Range("H100").Value = 20000618
$ws.Range("J100").Value = 696
$ws.Range("L100").Value = 1392
$ws.Range("N100").Value = -2474

$ws.Range("H107").Value = 532.7273
$ws.Range("I107").Value = 407.42856
$ws.Range("J107").Value = 752
$ws.Range("K107").Value = 1222.28568
$ws.Range("L107").Value = 2256
$ws.Range("M107").Value = 697.71432
$ws.Range("N107").Value = -6096

$ws.Range("H132").Value = 1856.2778
$ws.Range("I132").Value = 1941.6
$ws.Range("K132").Value = 5824.799999999999
$ws.Range("M132").Value = -3294.799999999999
